$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.066174745559692
$ws.Range("B1").Value = 2.693422079086304
$ws.Range("C1").Value = 2.899336814880371
$ws.Range("D1").Value = 3.727420806884766
$ws.Range("E1").Value = 5.116219043731689
